$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Items")

# --- Add new shared strings implicitly via cell values ---
# Row 2, 3, 4: add new zero-value cells in columns AS, AX, AY, BA
foreach ($r in 2..4) {
    $ws.Cells.Item($r, 45).Value = 0   # AS
    $ws.Cells.Item($r, 50).Value = 0   # AX
    $ws.Cells.Item($r, 51).Value = 0   # AY
    $ws.Cells.Item($r, 53).Value = 0   # BA
}

# --- Row 6: update id, and add several new zero cells ---
$ws.Cells.Item(6, 1).Value = 3563363            # A6 id
$ws.Cells.Item(6, 15).Value = 0                 # O6 copper_coin_cost
$ws.Cells.Item(6, 25).Value = 0                 # Y6 agi_mod
$ws.Cells.Item(6, 26).Value = 0                 # Z6 focus_mod
$ws.Cells.Item(6, 35).Value = 0                 # AI6 base_damage_mod_bonus
$ws.Cells.Item(6, 36).Value = 0                 # AJ6 base_healing_mod_bonus
$ws.Cells.Item(6, 37).Value = 0                 # AK6 base_ac_mod_bonus
$ws.Cells.Item(6, 38).Value = 0                 # AL6 fight_time_out_mod_bonus
$ws.Cells.Item(6, 39).Value = 0                 # AM6 move_time_out_mod_bonus
$ws.Cells.Item(6, 45).Value = 0                 # AS6 kingdom_damage
$ws.Cells.Item(6, 48).Value = 0                 # AV6 increase_stat_by
$ws.Cells.Item(6, 50).Value = 0                 # AX6 increase_skill_bonus_by
$ws.Cells.Item(6, 51).Value = 0                 # AY6 increase_skill_training_bonus_by
$ws.Cells.Item(6, 53).Value = 0                 # BA6 resurrection_chance
$ws.Cells.Item(6, 54).Value = 0                 # BB6 spell_evasion
$ws.Cells.Item(6, 55).Value = 0                 # BC6 artifact_annulment
$ws.Cells.Item(6, 56).Value = 0                 # BD6 healing_reduction
$ws.Cells.Item(6, 57).Value = 0                 # BE6 affix_damage_reduction
$ws.Cells.Item(6, 58).Value = 0                 # BF6 devouring_light
$ws.Cells.Item(6, 59).Value = 0                 # BG6 devouring_darkness

# --- Row 7: brand new item row ---
$ws.Cells.Item(7, 1).Value = 4418294            # A7 id
$ws.Cells.Item(7, 3).Value = "Blacksmiths Anvil" # C7 name
$ws.Cells.Item(7, 4).Value = "artifact"          # D7 type
$ws.Cells.Item(7, 7).Value = "A small light weight anvil that glistens with magic. A blacksmith might like this!" # G7 description
$ws.Cells.Item(7, 17).Value = 0.45              # Q7 base_damage_mod
$ws.Cells.Item(7, 18).Value = 0.1               # R7 base_healing_mod
$ws.Cells.Item(7, 19).Value = 0.3               # S7 base_ac_mod
$ws.Cells.Item(7, 20).Value = 0.4               # T7 str_mod
$ws.Cells.Item(7, 21).Value = 0.15              # U7 dur_mod
$ws.Cells.Item(7, 22).Value = 0.1               # V7 dex_mod
$ws.Cells.Item(7, 29).Value = 1                 # AC7 can_drop
$ws.Cells.Item(7, 48).Value = 0                 # AV7 increase_stat_by
$ws.Cells.Item(7, 65).Value = 0                 # BM7 holy_stacks
$ws.Cells.Item(7, 66).Value = 0                 # BN7 ambush_chance
$ws.Cells.Item(7, 67).Value = 0                 # BO7 ambush_resistance
$ws.Cells.Item(7, 68).Value = 0                 # BP7 counter_chance
$ws.Cells.Item(7, 69).Value = 0                 # BQ7 counter_resistance
$ws.Cells.Item(7, 72).Value = "Enraged Muscles"  # BT7 item_skill_id
